$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: set A2 to "CN"
$ws.Range("A2").Value = "CN"

# Row 3: fill in the new Stickler Thief entry
$ws.Range("A3").Value = "GFG"
$ws.Range("B3").Value = "GFG/Coding ninja"
$ws.Range("C3").Value = "Stickler Thief/Maximum sum of non-adjacent elements"
$ws.Range("D3").Value = "Java"
$ws.Range("E3").Value = "DP(Recurrsion+Memonization+Tabulation+space optimization)"

# E3 should keep the wrapped / top-left alignment used elsewhere in column E
$ws.Range("E3").HorizontalAlignment = -4131
$ws.Range("E3").VerticalAlignment = -4160
$ws.Range("E3").WrapText = $true

# Row now holds two lines of wrapped text, same as row 1
$ws.Rows(3).RowHeight = 30

# Update selection to A3
$ws.Range("A3").Select()
